$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("General")
$ws.Range("B5").Value = "$"
$ws.Range("B2").Value = 0.09
